$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewLoanInput")
$ws.Activate()

# Rename/replace the old, one-off scenario name in B2 with the new grouped
# scenario identifier (scenario naming convention grouping related scenarios).
$ws.Range("B2").Value = "3500-RBI-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME-DISBURSE-FEE-%APR-AMT-Reg-PERIODIC"

# Align B2's cell formatting with the plain fill style already used elsewhere
# on the sheet (e.g. B8), since the old bespoke "fill + alignment" style is
# no longer needed once B2 stops being a header-ish wrapped cell.
$ws.Range("B8").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to B2
$ws.Range("B2").Select()
